$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in operation import schema header: "TPJ" -> "TJ"
$ws.Range("I1").Value = "TJ"

# Reflect the last-selected cell as recorded in the saved file
$ws.Range("I1").Select()
